# Applies the cryptocurrency price/volume refresh described in the commit
# "Updated cryptos list on Sat Jul 29 19:14:35 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.345.87"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.98"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7116"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.21"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08092"
$ws.Range("E8").Value = "  +4.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3126"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.24"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08390"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.67"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.252"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7194"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.41"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.231"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008385"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.334.75"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.15"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.25"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.131.60"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.801"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1595"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.92"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.062"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.507"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.426"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.351"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  -7.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05374"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.949"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7515"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.178"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01890"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.288.88"
$ws.Range("E39").Value = "  +11.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.737"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.590"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.01"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8921"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.28"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  +8.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.021.53"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.802"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5207"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.475"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4366"
$ws.Range("E51").Value = "  +1.52%  "
